$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("(u6605935)", $true, $false, $false, $false, $false,
                $true, 1, $false, "", 0)
$r.InsertBefore(" ")
